$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("B2").Value = "Jesse Hare"
$ws.Range("G2").Value = 11

# Row 4 - Final Implementation (Individual)
$ws.Range("A4").Value = "Final Implementation"
$ws.Range("C4").Value = "I"
$ws.Range("D4").Value = "10/7/2019"
$ws.Range("E4").Value = 0.375
$ws.Range("F4").Value = 0.58333333333333337
$ws.Range("H4").Value = 5

# Row 5 - Unit Testing (Group)
$ws.Range("A5").Value = "Unit Testing"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = "10/8/2019"
$ws.Range("E5").Value = 0.375
$ws.Range("F5").Value = 0.58333333333333337
$ws.Range("G5").Value = 5

# Row 6 - Documentation (Group)
$ws.Range("A6").Value = "Documentation"
$ws.Range("C6").Value = "G"
$ws.Range("D6").Value = "10/9/2019"
$ws.Range("E6").Value = 0.375
$ws.Range("F6").Value = 0.58333333333333337
$ws.Range("G6").Value = 5

# Row 7 - Integration Testing (Individual)
$ws.Range("A7").Value = "Integration Testing"
$ws.Range("C7").Value = "I"
$ws.Range("D7").Value = "10/10/2019"
$ws.Range("E7").Value = 0.375
$ws.Range("F7").Value = 0.58333333333333337
$ws.Range("H7").Value = 5

$ws.Range("F4:F7").Select()
